{"js": "// Add a new paragraph \"Test stuuf\" after the last paragraph of the document\n// (the \"Peachy matey\" paragraph, which currently carries the trailing\n// _GoBack bookmark). When a user types new content at the end of a Word\n// document, Word automatically relocates the _GoBack bookmark onto the\n// newly edited paragraph, so mirror that by moving the bookmark too.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\"Test stuuf\", \"After\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\nnewParagraph.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$end = $d.Content\n$end.Collapse(0)   # wdCollapseEnd (0 = wdCollapseEnd)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.InsertAfter(\"Test stuuf\")\n"}
